$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Cel" column (C) to hold the
# new "Start" (departure airport) attribute.
$ws.Columns("C:C").Insert()

# --- Header row ---
$ws.Range("C1").Value = "Start"
$ws.Range("F1").Value = "Datum"

# --- Start (departure) column ---
$ws.Range("C2").Value = "DEB"
$ws.Range("C3").Value = "BUD"
$ws.Range("C4").Value = "BUD"
$ws.Range("C5").Value = "DEB"

# --- Datum (date/time) column ---
$ws.Range("F2").Value = 45626.5
$ws.Range("F3").Value = 45672.25
$ws.Range("F4").Value = 45650.833333333336
$ws.Range("F5").Value = 45626.416666666664
$ws.Range("F2:F5").NumberFormat = "yyyy/ m/ d. h:mm;@"

# Widen the new Datum column to fit its contents, like Excel's AutoFit.
$ws.Columns("F:F").ColumnWidth = 16.7109375

# Page setup (portrait, A4) as seen by the author when printing the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to the next empty row, below the new column.
$ws.Range("F6").Select()

# Reposition the workbook window, matching the author's on-screen layout.
$excel.ActiveWindow.Left = 4230
$excel.ActiveWindow.Top = 2730
